$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '58.284.46'
$ws.Range("E2").Value = '  +2.71%  '

$ws.Range("D3").Value = '2.370.16'
$ws.Range("E3").Value = '  +1.27%  '

$ws.Range("E4").Value = '  -0.03%  '

$ws.Range("D5").Value = '543.93'
$ws.Range("E5").Value = '  +5.67%  '

$ws.Range("D6").Value = '135.36'
$ws.Range("E6").Value = '  +1.09%  '

$ws.Range("D7").Value = '1.00'
$ws.Range("E7").Value = '  +0.12%  '

$ws.Range("D8").Value = '0.537'
$ws.Range("E8").Value = '  +0.58%  '

$ws.Range("D9").Value = '2.369.83'
$ws.Range("E9").Value = '  +0.90%  '

$ws.Range("E10").Value = '  +1.32%  '

$ws.Range("D11").Value = '5.44'
$ws.Range("E11").Value = '  +2.16%  '

$ws.Range("D12").Value = '0.154'
$ws.Range("E12").Value = '  +0.75%  '

$ws.Range("D13").Value = '0.354'
$ws.Range("E13").Value = '  +4.23%  '

$ws.Range("B14").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C14").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D14").Value = '2.788.17'
$ws.Range("E14").Value = '  +1.03%  '

$ws.Range("B15").Value = 'Avalanche'
$ws.Range("C15").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D15").Value = '23.74'
$ws.Range("E15").Value = '  -0.57%  '

$ws.Range("D16").Value = '58.165.74'
$ws.Range("E16").Value = '  +2.58%  '

$ws.Range("E17").Value = '  +0.37%  '

$ws.Range("D18").Value = '2.357.79'
$ws.Range("E18").Value = '  -0.16%  '

$ws.Range("D19").Value = '339.93'
$ws.Range("E19").Value = '  +4.18%  '

$ws.Range("E20").Value = '  +1.00%  '

$ws.Range("D21").Value = '4.23'
$ws.Range("E21").Value = '  +0.90%  '

$ws.Range("D22").Value = '6.90'
$ws.Range("E22").Value = '  +3.16%  '

$ws.Range("D23").Value = '0.999'
$ws.Range("E23").Value = '  -0.08%  '

$ws.Range("D24").Value = '62.33'
$ws.Range("E24").Value = '  +1.84%  '

$ws.Range("E25").Value = '  +2.54%  '

$ws.Range("D26").Value = '8.58'
$ws.Range("E26").Value = '  -1.30%  '

$ws.Range("D27").Value = '0.996'
$ws.Range("E27").Value = '  -0.21%  '

$ws.Range("E28").Value = '  +6.82%  '

$ws.Range("D29").Value = '175.08'
$ws.Range("E29").Value = '  +3.81%  '

$ws.Range("E30").Value = '  +4.85%  '

$ws.Range("E31").Value = '  +1.69%  '

$ws.Range("E32").Value = '  +0.11%  '

$ws.Range("D33").Value = '18.60'
$ws.Range("E33").Value = '  +0.97%  '

$ws.Range("E34").Value = '  +13.09%  '

$ws.Range("D35").Value = '0.999'
$ws.Range("E35").Value = '  -0.01%  '

$ws.Range("D36").Value = '0.999'
$ws.Range("E36").Value = '  +0.24%  '

$ws.Range("D37").Value = '1.26'
$ws.Range("E37").Value = '  -0.91%  '

$ws.Range("D38").Value = '4.11'
$ws.Range("E38").Value = '  +3.04%  '

$ws.Range("D39").Value = '1.62'
$ws.Range("E39").Value = '  +3.20%  '

$ws.Range("D40").Value = '39.43'
$ws.Range("E40").Value = '  +2.10%  '

$ws.Range("D41").Value = '149.92'
$ws.Range("E41").Value = '  -0.60%  '

$ws.Range("E42").Value = '  +0.57%  '

$ws.Range("E43").Value = '  +1.35%  '

$ws.Range("D44").Value = '285.16'
$ws.Range("E44").Value = '  +1.22%  '

$ws.Range("D45").Value = '0.0930'
$ws.Range("E45").Value = '  +0.43%  '

$ws.Range("D46").Value = '0.0505'
$ws.Range("E46").Value = '  +0.89%  '

$ws.Range("D47").Value = '18.96'
$ws.Range("E47").Value = '  +3.20%  '

$ws.Range("E48").Value = '  +0.90%  '

$ws.Range("D49").Value = '0.0219'
$ws.Range("E49").Value = '  +1.24%  '

$ws.Range("D50").Value = '17.60'
$ws.Range("E50").Value = '  +2.79%  '

$ws.Range("D51").Value = '0.383'
$ws.Range("E51").Value = '  +8.74%  '
